$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row at row 10 for "EXOSIRYLIC 500 MG 20 F.C.TABS."
#    (pushes FORTAZEDIM .. مناديل سولو سحب صغيره, totals row and footer down
#    by one row).
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).Insert()
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "EXOSIRYLIC 500 MG 20 F.C.TABS."
$ws.Range("H10").Value = "0:0"
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "1"
$ws.Range("L10").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"
$ws.Range("N10").Value = "194.00"
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "194.0000"
$ws.Range("P10").NumberFormat = "0.00"
$ws.Range("Q10").Value = "1:0"

# ---------------------------------------------------------------------------
# 2) Insert a new data row at row 13 for "PANADOL ADVANCE 500 MG 48 TABLETS"
#    (after MAVILOR, which is now row 12; pushes SPASMOFEN .. footer further
#    down by one more row).
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Range("A12:Q12").Copy()
$ws.Range("A13:Q13").PasteSpecial(-4122)
$ws.Range("A13:B13").Merge()
$ws.Range("C13:G13").Merge()
$ws.Range("H13:K13").Merge()
$ws.Range("L13:M13").Merge()
$ws.Range("N13:O13").Merge()

$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "PANADOL ADVANCE 500 MG 48 TABLETS"
$ws.Range("H13").Value = "1:2"
$ws.Range("L13").NumberFormat = "@"
$ws.Range("L13").Value = "1"
$ws.Range("L13").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"
$ws.Range("N13").Value = "92.00"
$ws.Range("P13").NumberFormat = "@"
$ws.Range("P13").Value = "92.0000"
$ws.Range("P13").NumberFormat = "0.00"
$ws.Range("Q13").Value = "1:0"

# ---------------------------------------------------------------------------
# 3) Renumber the remaining rows (they were copied along with their old
#    sequence numbers in column A when the rows shifted down).
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = 5   # FORTAZEDIM
$ws.Range("A12").Value = 6   # MAVILOR
$ws.Range("A14").Value = 8   # SPASMOFEN
$ws.Range("A15").Value = 9   # سرنجات 3 سم
$ws.Range("A16").Value = 10  # سرنجات 5 سم
$ws.Range("A17").Value = 11  # كالونا
$ws.Range("A18").Value = 12  # محلول رينجر
$ws.Range("A19").Value = 13  # مناديل سولو سحب صغيره

# ---------------------------------------------------------------------------
# 4) Update the total shown in the totals row (now row 20): the workbook adds
#    the two new product prices (194.00 + 92.00) to the previous total.
# ---------------------------------------------------------------------------
$ws.Range("P20").Value = 942.76499999999999

# ---------------------------------------------------------------------------
# 5) Update the generated-on timestamp in the footer (now row 21).
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Friday, 20 June, 2025 5:24 PM"
